$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 34 (shifts rows 34-49 down to 35-50)
$ws.Rows.Item(34).Insert()

# Update row 33's title (column B) to add "for Social Geography"
$ws.Cells.Item(33, 2).Value = "A course on Spatial Data Science for Social Geography"

# Fill the new row 34 with the micro course entry
$ws.Cells.Item(34, 1).Value = "Charles University in Prague"
$ws.Cells.Item(34, 2).Value = "A course on Spatial Data Science"
$ws.Cells.Item(34, 3).Value = "https://martinfleischmann.net/sds/micro/"
$ws.Cells.Item(34, 4).Value = "GeoData"

# Apply the same style as other data cells (wrap text) to the new row
$ws.Range("A34:D34").WrapText = $true
$ws.Rows.Item(34).AutoFit()

# Update view state to match (scroll position / active selection)
$ws.Application.ActiveWindow.ScrollRow = 26
$ws.Range("A33").Select()
